$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows from bottom to top so row numbers don't shift unexpectedly
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(2).Delete()
